$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7456961274147034
$ws.Range("B1").Value = 1.045445561408997
$ws.Range("C1").Value = 1.415834665298462
$ws.Range("D1").Value = 4.509974002838135
$ws.Range("E1").Value = 2.34572172164917
